# Update input data values on the "output" sheet (rows 2-25, various columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 109
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.675
$ws.Range("M2").Value = -100
$ws.Range("O2").Value = 0
$ws.Range("Q2").Value = 60
$ws.Range("T2").Value = 50
$ws.Range("W2").Value = -45

$ws.Range("J3").Value = 15
$ws.Range("L3").Value = 0.675
$ws.Range("M3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("Q3").Value = 60
$ws.Range("T3").Value = 50
$ws.Range("W3").Value = -50

$ws.Range("J4").Value = 87.4211
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = -68.4211
$ws.Range("O4").Value = 0
$ws.Range("Q4").Value = 60
$ws.Range("T4").Value = 50
$ws.Range("W4").Value = -55

$ws.Range("J5").Value = 11
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("Q5").Value = 60
$ws.Range("T5").Value = 50
$ws.Range("W5").Value = -55

$ws.Range("J6").Value = 21
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("Q6").Value = 60
$ws.Range("T6").Value = 50
$ws.Range("W6").Value = -50

$ws.Range("B7").Value = 190
$ws.Range("J7").Value = 30
$ws.Range("M7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("Q7").Value = 60
$ws.Range("T7").Value = 50
$ws.Range("W7").Value = -40

$ws.Range("J8").Value = 110
$ws.Range("K8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0

$ws.Range("J9").Value = 84
$ws.Range("K9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0

$ws.Range("J10").Value = 83
$ws.Range("K10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0

$ws.Range("J11").Value = 81
$ws.Range("K11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0

$ws.Range("B12").Value = 300
$ws.Range("J12").Value = 94
$ws.Range("K12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0

$ws.Range("B13").Value = 320
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0.557895
$ws.Range("N13").Value = 84
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0

$ws.Range("B14").Value = 280
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0.273684
$ws.Range("N14").Value = 54
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0

$ws.Range("B15").Value = 260
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0.115789
$ws.Range("N15").Value = 30
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0

$ws.Range("K16").Value = 0
$ws.Range("N16").Value = 22
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0

$ws.Range("B17").Value = 200
$ws.Range("K17").Value = -10.9197
$ws.Range("L17").Value = 0.0526316
$ws.Range("M17").Value = -11.0803
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0

$ws.Range("B18").Value = 180
$ws.Range("K18").Value = -36
$ws.Range("L18").Value = 0.0526316
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0

$ws.Range("B19").Value = 190
$ws.Range("K19").Value = 0
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0

$ws.Range("B20").Value = 240
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0

$ws.Range("B21").Value = 280
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0

$ws.Range("B22").Value = 325
$ws.Range("J22").Value = 165
$ws.Range("K22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0

$ws.Range("B23").Value = 350
$ws.Range("J23").Value = 190
$ws.Range("K23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0

$ws.Range("B24").Value = 300
$ws.Range("J24").Value = 140
$ws.Range("K24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0

$ws.Range("B25").Value = 250
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
